# Update the build-version timestamp embedded in the "About" sheet and in the
# "Boundaries and methane sources" sheet's build_version column.
#
# Old version string: mines - January 30 (built on January 30 2026 16.19.47 EST)
# New version string: mines - January 30 (built on February 02 2026 12.49.33 EST)

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- "About" sheet -------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# A2: "Version: mines - January 30 (built on January 30 2026 16.19.47 EST)"
$cellA2 = $wsAbout.Range("A2")
$textA2 = $cellA2.Value()
$cellA2.Value = $textA2.Replace($oldVersion, $newVersion)

# A6: Recommended Citation text containing the version string
$cellA6 = $wsAbout.Range("A6")
$textA6 = $cellA6.Value()
$cellA6.Value = $textA6.Replace($oldVersion, $newVersion)

# --- "Boundaries and methane sources" sheet -------------------------------
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S ("build_version") holds the same version string for every data row.
$usedRange = $wsData.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S
    $cellText = $cell.Value()
    if ($cellText -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
